$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-302 (years 1600-1900): set column B to 0
for ($r = 2; $r -le 302; $r++) {
    $ws.Cells.Item($r, 2).Value = 0
}

# Rows 303-452 (years 1901-2050): set column B to new computed values
$ws.Cells.Item(303, 2).Value = 29168971.06642462
$ws.Cells.Item(304, 2).Value = 58791219.39275392
$ws.Cells.Item(305, 2).Value = 88874426.49755873
$ws.Cells.Item(306, 2).Value = 119426371.0020093
$ws.Cells.Item(307, 2).Value = 150454926.2052194
$ws.Cells.Item(308, 2).Value = 181968057.3532576
$ws.Cells.Item(309, 2).Value = 213973818.5898039
$ws.Cells.Item(310, 2).Value = 246480349.576483
$ws.Cells.Item(311, 2).Value = 279495871.771022
$ws.Cells.Item(312, 2).Value = 313028684.3516133
$ws.Cells.Item(313, 2).Value = 347087159.7760499
$ws.Cells.Item(314, 2).Value = 381679738.9644697
$ws.Cells.Item(315, 2).Value = 416814926.0949376
$ws.Cells.Item(316, 2).Value = 452501283.0014567
$ws.Cells.Item(317, 2).Value = 488747423.1644644
$ws.Cells.Item(318, 2).Value = 525562005.2844399
$ws.Cells.Item(319, 2).Value = 562953726.4297929
$ws.Cells.Item(320, 2).Value = 600931314.7509103
$ws.Cells.Item(321, 2).Value = 639503521.753009
$ws.Cells.Item(322, 2).Value = 678679114.1211914
$ws.Cells.Item(323, 2).Value = 718466865.0920875
$ws.Cells.Item(324, 2).Value = 758875545.3673648
$ws.Cells.Item(325, 2).Value = 799913913.5655406
$ws.Cells.Item(326, 2).Value = 841590706.2095883
$ws.Cells.Item(327, 2).Value = 883914627.2491076
$ws.Cells.Item(328, 2).Value = 926894337.1171376
$ws.Cells.Item(329, 2).Value = 970538441.3230584
$ws.Cells.Item(330, 2).Value = 1014855478.584517
$ws.Cells.Item(331, 2).Value = 1059853908.50287
$ws.Cells.Item(332, 2).Value = 1105542098.788281
$ws.Cells.Item(333, 2).Value = 1151928312.042248
$ws.Cells.Item(334, 2).Value = 1199020692.10724
$ws.Cells.Item(335, 2).Value = 1246827249.994843
$ws.Cells.Item(336, 2).Value = 1295355849.405805
$ws.Cells.Item(337, 2).Value = 1344614191.857341
$ws.Cells.Item(338, 2).Value = 1394609801.435022
$ws.Cells.Item(339, 2).Value = 1445350009.188687
$ws.Cells.Item(340, 2).Value = 1496841937.193955
$ws.Cells.Item(341, 2).Value = 1549092482.302951
$ws.Cells.Item(342, 2).Value = 1602108299.610108
$ws.Cells.Item(343, 2).Value = 1655895785.661032
$ws.Cells.Item(344, 2).Value = 1710461061.434589
$ws.Cells.Item(345, 2).Value = 1765809955.130544
$ws.Cells.Item(346, 2).Value = 1821947984.797151
$ws.Cells.Item(347, 2).Value = 1878880340.835336
$ws.Cells.Item(348, 2).Value = 1936611868.417949
$ws.Cells.Item(349, 2).Value = 1995147049.864776
$ws.Cells.Item(350, 2).Value = 2054489987.015719
$ws.Cells.Item(351, 2).Value = 2114644383.646532
$ws.Cells.Item(352, 2).Value = 2175613527.973079
$ws.Cells.Item(353, 2).Value = 2237400275.291818
$ws.Cells.Item(354, 2).Value = 2300007030.805539
$ws.Cells.Item(355, 2).Value = 2363435732.684841
$ws.Cells.Item(356, 2).Value = 2427687835.416883
$ws.Cells.Item(357, 2).Value = 2492764293.493985
$ws.Cells.Item(358, 2).Value = 2570363238.447273
$ws.Cells.Item(359, 2).Value = 2648938285.700133
$ws.Cells.Item(360, 2).Value = 2728489061.363492
$ws.Cells.Item(361, 2).Value = 2809014551.675253
$ws.Cells.Item(362, 2).Value = 2890513087.532095
$ws.Cells.Item(363, 2).Value = 2972982329.814586
$ws.Cells.Item(364, 2).Value = 3056419255.567338
$ws.Cells.Item(365, 2).Value = 3140820145.094602
$ws.Cells.Item(366, 2).Value = 3226180570.03037
$ws.Cells.Item(367, 2).Value = 3312495382.440287
$ws.Cells.Item(368, 2).Value = 3399758705.010432
$ws.Cells.Item(369, 2).Value = 3487963922.375803
$ws.Cells.Item(370, 2).Value = 3577103673.638399
$ws.Cells.Item(371, 2).Value = 3667169846.121826
$ws.Cells.Item(372, 2).Value = 3758153570.406016
$ws.Cells.Item(373, 2).Value = 3818205420.12492
$ws.Cells.Item(374, 2).Value = 3878820230.967291
$ws.Cells.Item(375, 2).Value = 3939989133.00454
$ws.Cells.Item(376, 2).Value = 4001702684.972258
$ws.Cells.Item(377, 2).Value = 4063950873.241882
$ws.Cells.Item(378, 2).Value = 4126723111.715204
$ws.Cells.Item(379, 2).Value = 4190008242.652295
$ws.Cells.Item(380, 2).Value = 4253794538.440477
$ws.Cells.Item(381, 2).Value = 4318069704.308847
$ws.Cells.Item(382, 2).Value = 4382820881.989694
$ws.Cells.Item(383, 2).Value = 4423565963.521686
$ws.Cells.Item(384, 2).Value = 4464547742.524136
$ws.Cells.Item(385, 2).Value = 4505754482.243301
$ws.Cells.Item(386, 2).Value = 4547174000.495899
$ws.Cells.Item(387, 2).Value = 4588793671.96349
$ws.Cells.Item(388, 2).Value = 4630600431.025314
$ws.Cells.Item(389, 2).Value = 4672580775.122935
$ws.Cells.Item(390, 2).Value = 4714720768.649725
$ws.Cells.Item(391, 2).Value = 4757006047.357859
$ws.Cells.Item(392, 2).Value = 4799421823.275569
$ws.Cells.Item(393, 2).Value = 4874265749.425773
$ws.Cells.Item(394, 2).Value = 4949430133.473125
$ws.Cells.Item(395, 2).Value = 5024894800.554424
$ws.Cells.Item(396, 2).Value = 5100639155.528324
$ws.Cells.Item(397, 2).Value = 5176642201.600743
$ws.Cells.Item(398, 2).Value = 5252882559.695244
$ws.Cells.Item(399, 2).Value = 5329338488.542479
$ws.Cells.Item(400, 2).Value = 5405987905.463942
$ws.Cells.Item(401, 2).Value = 5482808407.827065
$ws.Cells.Item(402, 2).Value = 5559777295.15028
$ws.Cells.Item(403, 2).Value = 5628875646.307813
$ws.Cells.Item(404, 2).Value = 5698034308.624913
$ws.Cells.Item(405, 2).Value = 5767230716.052298
$ws.Cells.Item(406, 2).Value = 5836442075.570866
$ws.Cells.Item(407, 2).Value = 5905645389.333227
$ws.Cells.Item(408, 2).Value = 5974817477.390794
$ws.Cells.Item(409, 2).Value = 6043935001.013852
$ws.Cells.Item(410, 2).Value = 6112974486.61493
$ws.Cells.Item(411, 2).Value = 6181912350.289228
$ws.Cells.Item(412, 2).Value = 6250724922.988564
$ws.Cells.Item(413, 2).Value = 6319088403.426763
$ws.Cells.Item(414, 2).Value = 6387277870.807488
$ws.Cells.Item(415, 2).Value = 6455269576.881297
$ws.Cells.Item(416, 2).Value = 6523039793.809959
$ws.Cells.Item(417, 2).Value = 6590564842.732728
$ws.Cells.Item(418, 2).Value = 6657821123.101306
$ws.Cells.Item(419, 2).Value = 6724785142.804435
$ws.Cells.Item(420, 2).Value = 6791433549.100411
$ws.Cells.Item(421, 2).Value = 6857743160.371743
$ws.Cells.Item(422, 2).Value = 6923690998.71082
$ws.Cells.Item(423, 2).Value = 6989254323.339246
$ws.Cells.Item(424, 2).Value = 7054410664.855399
$ws.Cells.Item(425, 2).Value = 7119137860.29592
$ws.Cells.Item(426, 2).Value = 7183414088.986096
$ws.Cells.Item(427, 2).Value = 7247217909.142545
$ws.Cells.Item(428, 2).Value = 7310528295.178492
$ws.Cells.Item(429, 2).Value = 7373324675.647616
$ws.Cells.Item(430, 2).Value = 7435586971.747218
$ws.Cells.Item(431, 2).Value = 7497295636.285242
$ws.Cells.Item(432, 2).Value = 7558431692.99842
$ws.Cells.Item(433, 2).Value = 7618976776.091414
$ws.Cells.Item(434, 2).Value = 7678913169.84861
$ws.Cells.Item(435, 2).Value = 7738223848.152268
$ws.Cells.Item(436, 2).Value = 7796892513.722472
$ws.Cells.Item(437, 2).Value = 7854903636.877271
$ws.Cells.Item(438, 2).Value = 7912242493.594046
$ws.Cells.Item(439, 2).Value = 7968895202.638242
$ws.Cells.Item(440, 2).Value = 8024848761.510782
$ws.Cells.Item(441, 2).Value = 8080091080.953601
$ws.Cells.Item(442, 2).Value = 8134611017.742414
$ws.Cells.Item(443, 2).Value = 8188398405.488072
$ws.Cells.Item(444, 2).Value = 8241444083.163219
$ws.Cells.Item(445, 2).Value = 8293739921.069029
$ws.Cells.Item(446, 2).Value = 8345278843.958698
$ws.Cells.Item(447, 2).Value = 8396054851.039561
$ws.Cells.Item(448, 2).Value = 8446063032.584667
$ws.Cells.Item(449, 2).Value = 8495299582.898228
$ws.Cells.Item(450, 2).Value = 8543761809.395732
$ws.Cells.Item(451, 2).Value = 8591448137.581196
$ws.Cells.Item(452, 2).Value = 8638358111.728783
